$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45970, 0, 7.583),
    @(45970.01041666666, 0.12, 11.998),
    @(45970.02083333334, 1.23, 1.44),
    @(45970.03125, 3.624, 3.678),
    @(45970.04166666666, 3.204, 0.61),
    @(45970.05208333334, 2.952, 0.5600000000000001),
    @(45970.0625, 2.939, 1.127),
    @(45970.07291666666, 0.145, 3.377),
    @(45970.08333333334, 0, 43.665),
    @(45970.09375, 0.916, 9.532),
    @(45970.10416666666, 0.028, 5.968),
    @(45970.11458333334, 0, 22.399),
    @(45970.125, 0, 10.42),
    @(45970.13541666666, 0, 13.868),
    @(45970.14583333334, 0, 22.114),
    @(45970.15625, 0, 39.259),
    @(45970.16666666666, 0, 42.196),
    @(45970.17708333334, 0, 57.057),
    @(45970.1875, 0, 42.629),
    @(45970.19791666666, 0, 26.289),
    @(45970.20833333334, 0, 34.201),
    @(45970.21875, 0, 36.947),
    @(45970.22916666666, 0, 30.282),
    @(45970.23958333334, 0, 36.687),
    @(45970.25, 0, 49.486),
    @(45970.26041666666, 0, 45.473),
    @(45970.27083333334, 0, 24.09),
    @(45970.28125, 0.499, 4.425),
    @(45970.29166666666, 0, 11.828),
    @(45970.30208333334, 0.374, 4.843),
    @(45970.3125, 23.194, 0)
)

$row = 2
foreach ($line in $data) {
    $ws.Cells.Item($row, 1).Value = $line[0]
    $ws.Cells.Item($row, 2).Value = $line[1]
    $ws.Cells.Item($row, 3).Value = $line[2]
    $row++
}

$ws.Cells.Item(32, 1).NumberFormat = $ws.Cells.Item(31, 1).NumberFormat
